# Zowe architecture deck update:
#   1. Refresh the auto-date ("datetimeFigureOut") placeholder text that lives
#      on the slide master and every slide layout from 2/8/2019 -> 2/12/2019.
#   2. Remove the "JMON" box (Rectangle 65) and its connector
#      (Elbow Connector 62) from slide 1, per "Update the Zowe architecture
#      to remove JMON".

$p = $ppt.ActivePresentation

$newDate = "2/12/2019"

# --- 1. Update the date placeholder on the slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- Update the date placeholder on every slide layout ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Remove the JMON rectangle and its elbow connector from slide 1 ---
$slide = $p.Slides.Item(1)

for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "Elbow Connector 62" -or $sh.Name -eq "Rectangle 65") {
        $sh.Delete()
    }
}
